# Dataprovider - put invalid and valid data in separate sheets for scalability
#
# The original workbook has a single "loginSheet" that mixes the one valid
# login data-row together with several invalid login data-rows, and a single
# "pythonCode" sheet that mixes the one valid try-editor data-row with the
# one invalid try-editor data-row.  This script splits each of those sheets
# into a "valid*" and "invalid*" sheet so the data-providers scale better.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Split "loginSheet" into "validLoginSheet" + "invalidLoginSheet"
# ---------------------------------------------------------------------
$loginSheet = $wb.Worksheets.Item("loginSheet")

# Duplicate the sheet and place the copy immediately before the original;
# the copy becomes the "valid" sheet (it only needs the header + the single
# valid data row), the original becomes the "invalid" sheet (it keeps the
# header + all the invalid data rows).
$loginSheet.Copy($loginSheet)

$validLoginSheet = $wb.Worksheets.Item("loginSheet (2)")
$invalidLoginSheet = $wb.Worksheets.Item("loginSheet")

$validLoginSheet.Name = "validLoginSheet"
$invalidLoginSheet.Name = "invalidLoginSheet"

# validLoginSheet keeps only the header row (1) and the valid data row (2);
# everything below (the invalid rows) is removed.
$validLoginSheet.Range("A3:D8").EntireRow.Delete()

# invalidLoginSheet keeps the header row (1) and the invalid rows (old rows
# 3-8); the valid data row (old row 2) is removed and the hyperlink that
# pointed at it goes with it.
$invalidLoginSheet.Hyperlinks.Delete()
$invalidLoginSheet.Range("A2:D2").EntireRow.Delete()

# ---------------------------------------------------------------------
# 2. Split "pythonCode" into "validTryEditor" + "invalidTryEditor"
# ---------------------------------------------------------------------
$pythonCode = $wb.Worksheets.Item("pythonCode")

# Duplicate the sheet, placing the copy right after the original; the
# original becomes the "valid" sheet (header + the one valid try-editor
# row), the copy becomes the "invalid" sheet (header + the one invalid
# try-editor row).
$pythonCode.Copy($null, $pythonCode)

$validTryEditor = $wb.Worksheets.Item("pythonCode")
$invalidTryEditor = $wb.Worksheets.Item("pythonCode (2)")

$validTryEditor.Name = "validTryEditor"
$invalidTryEditor.Name = "invalidTryEditor"

# validTryEditor keeps the header row (1) and the valid row (2); the
# invalid row (old row 3) is removed.
$validTryEditor.Range("A3:B3").EntireRow.Delete()

# invalidTryEditor keeps the header row (1) and the invalid row (old row 3,
# which becomes row 2); the valid row (old row 2) is removed.
$invalidTryEditor.Range("A2:B2").EntireRow.Delete()

# ---------------------------------------------------------------------
# 3. Make "invalidLoginSheet" the active/selected sheet, like before
# ---------------------------------------------------------------------
$invalidLoginSheet.Activate()
